$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pretty-printed JSON text (replaces the old single-line questions = [...] string)
$text = @'
questions = [
    {
        "title": "Your colleague has upgraded a Magento site from version 2.2.2 to 2.3.3 and found the following error on the browser:\u201c(Magento\\Framework\\Config\\Dom\\ValidationException): Element 'referenceContainer', attribute 'before': The attribute 'before' is not allowed.\u201d How can you fix this error?",
        "ques_type": 2,
        "options": [
            "Search before attribute in referenceContainer element in phtml files and remove it. ",
            "Search before attribute in referenceContainer element in layout xml files and remove the 'before' attribute. ",
            "Search before attribute in referenceContainer element in layout xml files and turn the 'before' attribute to true.",
            "Search before attribute in referenceContainer element in phtml files and turn the 'before' attribute to true."
        ],
        "score": "Search before attribute in referenceContainer element in layout xml files and remove the 'before' attribute."
    },
    {
        "title": "A merchant asks you to display the message just below \u201cAdd to Cart\u201d, like in the following image. Your colleague has extended addtocart.phtml to Magento_Checkout directory in your custom theme.How can you accomplish this task effectively, with minimum code?",
        "ques_type": 2,
        "options": [
            "Copy catalog_product_view.xml layout file to your theme and edit the layout to move the  container under the addtocart block.",
            "Copy the &ltdiv data-bind=\"scope: 'messages'\"&gt with its knockout block from messages.phtml and paste under product-addtocart-button element.",
            "Copy catalog_product_view.xml layout file to your theme and edit the layout by pasting the message block under the addtocart block.",
            "Write JavaScript function for class=\"box-tocart\" to append the message block using jQuery."
        ],
        "score": "Copy the &ltdiv data-bind=\"scope: 'messages'\"&gt with its knockout block from messages.phtml and paste under product-addtocart-button element."
    },
    {
        "title": "A merchant asked you to develop a feature so that specific information about a product is not accessible by an anonymous or unauthorized user on the web API call. To do that, your colleague has written the following code in extension_attributes.xml file: &ltextension_attributes for=\"Magento\\Catalog\\Api\\Data\\ProductInterface\"&gt\n       &ltattribute code=\"stock_item\" type=\"Magento\\CatalogInventory\\Api\\Data\\StockItemInterface\"&gt\n          &ltresources&gt &lt/resources&gt\n       &lt/attribute&gt\n&lt/extension_attributes&gt\nWhat modification can correct the above code, so it meets the requirements?",
        "ques_type": 2,
        "options": [
            "Write &ltresource ref=\"anonymous\"/&gt inside the &ltresources&gt tag.",
            "Write &ltresource ref=\"self\"/&gt inside the &ltresources&gt tag.",
            "Write &ltresource ref=\u201cMagento_CatalogInventory::cataloginventory\u201d/&gt inside the &ltresources&gt tag.",
            "Write &ltresource ref=\"api\"/&gt inside the &ltresources&gt tag."
        ],
        "score": "Write &ltresource ref=\u201cMagento_CatalogInventory::cataloginventory\u201d/&gt inside the &ltresources&gt tag."
    },
    {
        "title": "A merchant asks you to improve her site\u2019s loading time. What possible things can you do to make it faster?",
        "ques_type": 15,
        "options": [
            "Minify JavaScript and CSS files",
            "Disable Flat Categories and Products",
            "Make JavaScript and CSS internal",
            "Enable varnish cache",
            "Use only Magento Cache instead of varnish cache"
        ],
        "score": [
            "Minify JavaScript and CSS files",
            "Enable varnish cache"
        ]
    }
]
'@

# A1 currently holds a placeholder 0 with bold+bordered+centered formatting;
# A2 holds the real question text (shared string) with default formatting.
# Clear A1's formatting back to the workbook default, then overwrite its
# value with the updated text, then remove the now-duplicate row 2.
$ws.Range("A1").Style = "Normal"
$ws.Range("A1").Value = $text
$ws.Rows(1).EntireRow.AutoFit()
$ws.Range("A2").EntireRow.Delete()
